# "created 28 June folder, added TASK6 catalog file"
# Append a new Work-History row (row 18) to the Summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# --- New row of data -------------------------------------------------
$ws.Range("A18").Value = 44373
$ws.Range("B18").Value = "Reports on multiple dependency parameters (not less than 5 params)"
$ws.Range("C18").Value = "Bursting completed and results mail came"
$ws.Range("D18").Value = "No"

# --- Formatting to match the rest of the data rows --------------------
$ws.Range("A18").NumberFormat = "d-mmm"
$ws.Range("A18:D18").HorizontalAlignment = -4131
$ws.Range("A18:D18").VerticalAlignment = -4160

# --- Leave the workbook scrolled/selected where the author left off ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("B19").Select()
